# template_for_word_flicker.pptx:
#   "word and picture is put into the ppt according to the template for word flicker"
#
# The two placeholder-only slides that were used while building the template
# are removed (the deck goes back to just the master/layouts, ready to be
# used as a template), and the slide master's date placeholder is bumped to
# the day the template was finalized.

$p = $ppt.ActivePresentation

# Remove every slide from the deck (both slide1 and slide2), working from
# the end so indices stay valid as slides are deleted.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $p.Slides.Item($i).Delete()
}

# Update the slide master's date placeholder text.
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shape = $master.Shapes.Item($i)
    if ($shape.PlaceholderFormat.Type -eq 16) {
        $shape.TextFrame.TextRange.Text = "2018. 9. 20."
    }
}
